$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E6: Actual Start-date for the "Unit test" task -> 2019-09-16 (serial 43724)
# Copy the date format from an existing date cell (C3) so it reuses the
# same style (numFmtId 15 "d-mmm-yy") instead of minting a new number format.
$c3 = $ws.Range("C3")
$e6 = $ws.Range("E6")
$c3.Copy() | Out-Null
$e6.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$e6.Value = Get-Date -Year 2019 -Month 9 -Day 16 -Hour 0 -Minute 0 -Second 0

# G6: Status -> "In-Progress", shown in green font color (RGB 92D050)
$g6 = $ws.Range("G6")
$g6.Value = "In-Progress"
$g6.Font.Color = 5296274
